# M01 demonstration test plan: fill in the completed Course class test data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: developer name ---
$ws.Range("C3").Value = "Damien Altenburg"

# --- Row 7: __init__ / Attribute set to input value ---
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = 'name = "intermediate software development"
department = Department.COMPUTER_SCIENCE
credit_hours = 90'
$ws.Range("G7").Value = "Object initialized with the correct state"

# --- Row 8: __init__ / Exception raised when name is blank ---
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = 'name = ""
department = Department.COMPUTER_SCIENCE
credit_hours = 90'
$ws.Range("G8").Value = "ValueError: name cannot be blank."

# --- Row 9: __init__ / Exception raised when invalid department ---
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = 'name = "intermediate software development"
department = "computer science"
credit_hours = 90'
$ws.Range("G9").Value = "TypeError: department object must be a Department type."

# --- Row 10: __init__ / Exception raised when credit hours is not an int. ---
$ws.Range("D10").Value = "Exception raised when credit hours is not an int."
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = 'name = "intermediate software development"
department = Department.COMPUTER_SCIENCE
credit_hours = 34.5'
$ws.Range("G10").Value = "TypeError: credit_hours object must be an int type."

# --- Row 11: name / Returns name attribute ---
$ws.Range("D11").Value = "Returns name attribute"
$ws.Range("E11").Value = 'Object initialized.
name = "intermediate software development"
department = Department.COMPUTER_SCIENCE
credit_hours = 90'
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = '"intermediate software development"'

# --- Row 12: department / Returns department attribute ---
$ws.Range("D12").Value = "Returns department attribute"
$ws.Range("E12").Value = 'Object initialized.
name = "intermediate software development"
department = Department.COMPUTER_SCIENCE
credit_hours = 90'
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "Department.COMPUTER_SCIENCE"

# --- Row 13: credit_hours / Returns credit_hours attribute ---
$ws.Range("D13").Value = "Returns credit_hours attribute"
$ws.Range("E13").Value = 'Object initialized.
name = "intermediate software development"
department = Department.COMPUTER_SCIENCE
credit_hours = 90'
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = 90

# --- Row 14: __str__ / Returns string in expected format. ---
$ws.Range("D14").Value = "Returns string in expected format."
$ws.Range("E14").Value = 'Object initialized.
name = "intermediate software development"
department = Department.COMPUTER_SCIENCE
credit_hours = 90'
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = '"Course: Intermediate Software Development"
Department: Computer Science
Credit Hours: 90"'

# --- Row heights for the newly-filled rows (auto-sized by Excel after entry) ---
$ws.Rows.Item(8).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 105
$ws.Rows.Item(12).RowHeight = 105
$ws.Rows.Item(13).RowHeight = 105
$ws.Rows.Item(14).RowHeight = 105

# --- Blank template rows (15-26): normalize D:G to the plain (non-italic,
#     non-bold) style used by column C so the whole row reads consistently ---
$blankRows = 15..26
foreach ($r in $blankRows) {
    $rng = $ws.Range("D" + $r + ":G" + $r)
    $rng.Font.Italic = $false
    $rng.Font.Bold = $false
}

# --- Column widths (E and G got wider to fit the new content) ---
$ws.Columns.Item(5).ColumnWidth = 38.3
$ws.Columns.Item(7).ColumnWidth = 35.7

# --- Page setup: fit to one page wide, landscape, 57% scale, print area A1:G14 ---
$ws.PageSetup.Zoom = 57
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
$ws.PageSetup.Orientation = 2
$ws.PageSetup.PrintArea = '$A$1:$G$14'

# --- Selection reflects the print area ---
$ws.Range("A1:G14").Select()
